$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 19, mirroring the layout/style of row 18 (GFG / <text> / Java / date)
$ws.Range("A19").Value = "GFG"
$ws.Range("B19").Value = " Second Largest in Array"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "Java"
$ws.Range("D19").Value = 45000

$ws.Range("A19").HorizontalAlignment = $ws.Range("A18").HorizontalAlignment
$ws.Range("D19").NumberFormat = $ws.Range("D18").NumberFormat

# Update the active selection to match the new state
$ws.Range("L17").Select()
